$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: row 3
$ws.Range("M3").Value = "Interest Deduction Limitations"
# Row 4: Australia
$ws.Range("M4").Value = "1.5:1 debt-to-equity ratio (15:1 for financial institutions) applies"
# Row 5: Austria
$ws.Range("K5").Value = "Passive"
$ws.Range("L5").Value = "CFC with substantive economic activities exempted"
$ws.Range("M5").Value = "Informal 4:1 debt-to-equity ratio applies"
# Row 6: Belgium
$ws.Range("E6").Value = "None"
$ws.Range("K6").Value = "All Income Associated with Non-genuine Arrangements"
$ws.Range("L6").Value = "None"
$ws.Range("M6").Value = "Interest deductions limited to the higher of €3 million or 30% of EBITDA`n5:1 debt-to-equity ratio applies to intragroup loans`n1:1 debt-to-equity ratio applies to receivables from shareholders or directors, managers, and liquidators"
# Row 7: Canada
$ws.Range("E7").Value = "Countries with a tax treaty or Tax Information Exchange Agreement"
$ws.Range("L7").Value = "Multiple rules may exempt CFC from taxation"
$ws.Range("M7").Value = "1.5:1 debt-to-equity ratio applies"
# Row 8: Chile
$ws.Range("E8").Value = "N/A"
$ws.Range("L8").Value = "None"
$ws.Range("M8").Value = "3:1 debt-to-equity ratio applies`nA surtax for excessive-indebtedness can apply"
# Row 9: Czech Republic
$ws.Range("E9").Value = "EU member states and EEA member states"
$ws.Range("K9").Value = "Passive"
$ws.Range("L9").Value = "CFC with substantive economic activities exempted"
$ws.Range("M9").Value = "4:1 debt-to-equity ratio (6:1 debt-to-equity ratio for certain financial services companies) applies"
# Row 10: Denmark
$ws.Range("E10").Value = "EU member states and EEA member states or double tax treaty"
$ws.Range("K10").Value = "All Income"
$ws.Range("M10").Value = "4:1 debt-to-equity ratio applies`nInterest deductions are limited to 2.7% of assets`nInterest and depreciation deduction limited to 30% of EBITDA`nOther rules can apply"
# Row 11: Estonia
$ws.Range("E11").Value = "EU member states and EEA member states and Switzerland"
$ws.Range("K11").Value = "All Income Associated with Non-genuine Arrangements"
$ws.Range("M11").Value = "Interest deductions limited to the higher of €3 million or 30% of EBITDA"
# Row 12: Finland
$ws.Range("E12").Value = "EU member states and EEA member states"
$ws.Range("M12").Value = "Interest deductions limited to 25% of EBITDA`nNet interest expenses between non-related parties limited to €3 million"
# Row 13: France
$ws.Range("E13").Value = "Black-list countries are excluded"
$ws.Range("K13").Value = "All Income"
$ws.Range("M13").Value = "Interest deductions limited to the higher of €3 million or 30% of EBITDA"
# Row 14: Germany
$ws.Range("E14").Value = "None"
$ws.Range("M14").Value = "Interest deductions limited to the higher of €3 million or 30% of EBITDA"
# Row 15: Greece
$ws.Range("E15").Value = "EU member states"
$ws.Range("M15").Value = "Interest deductions limited to the higher of €3 million or 30% of EBITDA"
# Row 16: Hungary
$ws.Range("K16").Value = "Passive"
# Row 17: Iceland
$ws.Range("E17").Value = "None"
$ws.Range("K17").Value = "All Income"
$ws.Range("M17").Value = "Interest deductions limited to 30% of EBITDA`nRule does not apply if total interest paid does not exceed ISK 100 million `nOther exemptions can apply"
# Row 18: Ireland
$ws.Range("E18").Value = "EU member states and tax treaty countries"
$ws.Range("K18").Value = "All Income Associated with Non-genuine Arrangements"
$ws.Range("L18").Value = "CFC exempt if i) below certain profit and income thresholds; ii) transfer pricing rules apply; or iii) passes the essential purpose test."
$ws.Range("M18").Value = "None`nHowever, in specific cases, interest can be reclassified as a dividend"
# Row 19: Israel
$ws.Range("E19").Value = "N/A"
$ws.Range("K19").Value = "Passive"
$ws.Range("L19").Value = "None"
$ws.Range("M19").Value = "None"
# Row 20: Italy
$ws.Range("E20").Value = "Black-list countries are excluded"
# Row 21: Japan
$ws.Range("E21").Value = "None"
# Row 22: Korea
$ws.Range("E22").Value = "N/A"
$ws.Range("K22").Value = "All Income"
# Row 23: Latvia
$ws.Range("K23").Value = "All Income Associated with Non-genuine Arrangements"
# Row 24: Lithuania
$ws.Range("E24").Value = "Black-list countries are excluded"
$ws.Range("K24").Value = "Passive"
# Row 25: Luxembourg
$ws.Range("E25").Value = "None"
$ws.Range("K25").Value = "All Income Associated with Non-genuine Arrangements"
$ws.Range("L25").Value = "CFC exempt if i) not an artificial arrangement or ii) accounting profits  below €750,000 or less than 10% of operating costs"
# Row 26: Mexico
$ws.Range("E26").Value = "N/A"
$ws.Range("K26").Value = "All Income"
$ws.Range("L26").Value = "None"
# Row 28: New Zealand
$ws.Range("E28").Value = "None"
$ws.Range("K28").Value = "Passive"
# Row 29: Norway
$ws.Range("E29").Value = "EEA member states"
$ws.Range("L29").Value = "CFC exempt if located in EEA country and not an artificial arrangement or located in tax treaty country"
# Row 30: Poland
$ws.Range("E30").Value = "EU member states and EEA member states and Switzerland"
$ws.Range("L30").Value = "CFC exempt if located in EU or EEA and not an artificial arrangement"
$ws.Range("M30").Value = "Interest deductions limited to the higher of PLN 3 million or 30% of EBITDA"
# Row 31: Portugal
$ws.Range("E31").Value = "Black-list countries are excluded"
$ws.Range("K31").Value = "All Income"
$ws.Range("L31").Value = "CFC exempt if located in EU and EEA countries and not an artificial arrangement`nOther exemptions can apply"
$ws.Range("M31").Value = "Interest deductions limited to the higher of €1 million or 30% of EBITDA"
# Row 32: Slovak Republic
$ws.Range("E32").Value = "Tax treaty countries"
$ws.Range("K32").Value = "All Income Associated with Non-genuine Arrangements"
$ws.Range("L32").Value = "None"
# Row 33: Slovenia
$ws.Range("L33").Value = "Substantial economic activities exemption"
$ws.Range("M33").Value = "4:1 debt-to-equity ratio applies"
# Row 34: Spain
$ws.Range("E34").Value = "Black-list countries are excluded"
$ws.Range("K34").Value = "Passive"
$ws.Range("L34").Value = "CFC exempt if located in EU or EEA and not an artificial arrangement"
$ws.Range("M34").Value = "Interest deductions limited to the higher of €1 million or 30% of EBITDA"
# Row 35: Sweden
$ws.Range("J35").Value = "Yes"
$ws.Range("K35").Value = "All Income"
$ws.Range("L35").Value = "CFC exempt if located in EEA and not an artificial arrangement or located in white list countries"
$ws.Range("M35").Value = "Interest deductions limited to 30% of EBITDA"
# Row 36: Switzerland
$ws.Range("J36").Value = "No"
$ws.Range("K36").Value = "N/A"
$ws.Range("L36").Value = "N/A"
# Row 37: Turkey
$ws.Range("L37").Value = "None"
$ws.Range("M37").Value = "3:1 debt-to-equity ratio (6:1 for financial institutions) applies"
# Row 38: United Kingdom
$ws.Range("K38").Value = "All Income"
$ws.Range("L38").Value = "Various exemptions can apply"
$ws.Range("M38").Value = "Interest deductions limited to 30% of EBITDA"
# Row 39: United States
$ws.Range("E39").Value = "None"
$ws.Range("J39").Value = "Yes"
$ws.Range("K39").Value = "Passive"
$ws.Range("L39").Value = "Exemptions for foreign high-taxed income can apply"
$ws.Range("M39").Value = "Interest deductions limited to the sum of business interest income, 30% of adjusted taxable income, and floor plan financing interest"
